$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the current row 1245, shifting the existing
# rows 1245-1273 down to 1247-1275 (dimension grows from R1273 to R1275).
$ws.Rows.Item(1245).Insert()
$ws.Rows.Item(1245).Insert()

# New row 1245: weekly Ajo/Chino/Primera "$/caja 10 kilos" entry.
$ws.Cells.Item(1245,1).Value2  = 10
$ws.Cells.Item(1245,2).Value2  = "Vega Modelo de Temuco"
$ws.Cells.Item(1245,3).Value2  = "La Araucanía"
$ws.Cells.Item(1245,4).Value2  = 45239
$ws.Cells.Item(1245,5).Value2  = 9
$ws.Cells.Item(1245,6).Value2  = 100112003
$ws.Cells.Item(1245,7).Value2  = "Ajo"
$ws.Cells.Item(1245,8).Value2  = "Chino"
$ws.Cells.Item(1245,9).Value2  = "Primera"
$ws.Cells.Item(1245,10).Value2 = 450
$ws.Cells.Item(1245,11).Value2 = 23000
$ws.Cells.Item(1245,12).Value2 = 24000
$ws.Cells.Item(1245,13).Value2 = 23333
$ws.Cells.Item(1245,14).Value2 = "$/caja 10 kilos"
$ws.Cells.Item(1245,15).Value2 = "China"
$ws.Cells.Item(1245,16).Value2 = 2333
$ws.Cells.Item(1245,17).Value2 = 10
$ws.Cells.Item(1245,18).Value2 = "Hortaliza"

# New row 1246: weekly Ajo/Chino/Primera "$/malla 10 kilos" entry.
$ws.Cells.Item(1246,1).Value2  = 10
$ws.Cells.Item(1246,2).Value2  = "Vega Modelo de Temuco"
$ws.Cells.Item(1246,3).Value2  = "La Araucanía"
$ws.Cells.Item(1246,4).Value2  = 45239
$ws.Cells.Item(1246,5).Value2  = 9
$ws.Cells.Item(1246,6).Value2  = 100112003
$ws.Cells.Item(1246,7).Value2  = "Ajo"
$ws.Cells.Item(1246,8).Value2  = "Chino"
$ws.Cells.Item(1246,9).Value2  = "Primera"
$ws.Cells.Item(1246,10).Value2 = 200
$ws.Cells.Item(1246,11).Value2 = 26000
$ws.Cells.Item(1246,12).Value2 = 26000
$ws.Cells.Item(1246,13).Value2 = 26000
$ws.Cells.Item(1246,14).Value2 = "$/malla 10 kilos"
$ws.Cells.Item(1246,15).Value2 = "China"
$ws.Cells.Item(1246,16).Value2 = 2600
$ws.Cells.Item(1246,17).Value2 = 10
$ws.Cells.Item(1246,18).Value2 = "Hortaliza"
